# The "meta" sheet holds key/value metadata rows (tab/discontinuity,
# type/line, title/A discontinuity, y_lim/0,3, ...) followed by one blank
# row. This change adds a new "style" / "default" metadata row in the
# place of that blank row, and re-adds a blank row (with the same
# formatting as the other key cells) after it.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# Fill in the previously blank row 5 with the new style/default pair.
$meta.Range("A5").Value = "style"
$meta.Range("B5").Value = "default"

# A5 keeps the same "key" formatting as A1:A4 (bold, coloured text).
$meta.Range("A4").Copy()
$meta.Range("A5").PasteSpecial(-4122)

# Re-create the trailing blank row, now at row 6, with the same
# formatting the blank cell had at its old position (A5).
$meta.Range("A6").Value = $null
$meta.Range("A4").Copy()
$meta.Range("A6").PasteSpecial(-4122)

$excel.CutCopyMode = $false
